$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two schools ("Cal State Northridge", "Penn state") are no longer being
# tracked in this log, so drop those two rows entirely. Everything below
# shifts up by two rows (old row 24 "Acquafondata" lands on row 22).
$ws.Rows("2:3").Delete()

# Jot down the day/contact notes that were collected for a few schools.
$ws.Range("C5").Value = "Friday"
$ws.Range("C6").Value = "Monday/Tuesday - James"
$ws.Range("C18").Value = "Soon"

# Standardize the school names to their official names.
$ws.Cells.Replace("University of Nebraska", "University of Nebraska-Lincoln")
$ws.Cells.Replace("University of Nevada, Reno", "University of Nevada-Reno")
$ws.Cells.Replace("Hampden Sydney College", "Hampden-Sydney College")
$ws.Cells.Replace("Texas Christian", "Texas Christian University")
$ws.Cells.Replace("University of South Carolina", "University of South Carolina-Columbia")
$ws.Cells.Replace("USC", "University of Southern California")
$ws.Cells.Replace("Texas A &M", "Texas A & M University-College Station")
$ws.Cells.Replace("UC Riverside", "University of California-Riverside")
$ws.Cells.Replace("UC Santa Cruz", "University of California-Santa Cruz")
$ws.Cells.Replace("UC Irvine", "University of California-Irvine")
$ws.Cells.Replace("Cal State Fullerton", "California State University-Fullerton")
$ws.Cells.Replace("University of albany", "SUNY at Albany")

# Fill in the remaining "emailed"/"collected" tracking flags that changed.
$ws.Range("C3").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C19").Value = 0

# Highlight the schools that were renamed/standardized with a distinct font
# (Verdana 11, dark gray). Build the format once on the first cell, then copy
# it onto the rest so the style table doesn't grow a new font record per cell.
$styledRows = @(2, 3, 5, 9, 12, 14, 15, 16, 17, 18)
$firstStyled = $ws.Cells.Item($styledRows[0], 1)
$firstStyled.Font.Name = "Verdana"
$firstStyled.Font.Size = 11
$firstStyled.Font.Color = 3355443

$firstStyled.Copy()
foreach ($r in $styledRows[1..($styledRows.Length - 1)]) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Move the active selection.
$ws.Range("A10").Select()
